# feat: add 2022-Q1 data
#
# Before:  Sheet1 "2021-Q4"  |  Sheet2 "总计"
# After:   Sheet1 "2021-Q4"  |  Sheet2 "2022-Q1" (new)  |  Sheet3 "总计"
#
# Strategy:
#  - The existing "总计" worksheet (rId2 / sheetId 2) is repurposed in place to
#    become the new "2022-Q1" quarterly sheet (same physical sheet, renamed +
#    re-filled) so it keeps sheetId=2 / rId2, matching the target sheet order.
#  - A brand-new worksheet named "总计" is appended right after it, holding the
#    combined totals table (old totals row + the new 2022-Q1 row on top).
#  - xlPasteFormats is used to stamp the workbook's existing header/index cell
#    style (s="2": bold, centered, thin-bordered) onto newly created cells so
#    formatting matches the rest of the workbook exactly.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122
$xlLineStyleContinuous = 1
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlCenter = -4108
$xlTop = -4160

# ---------------------------------------------------------------------------
# 1) Repurpose the existing "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# Grab a cell that already carries the shared "header" style (bold, centered,
# thin border) so we can stamp it onto the new cells we are about to create.
$styleSource = $q1.Range("B1")

# -- Header row (B1:H1) ------------------------------------------------------
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

$styleSource.Copy()
$q1.Range("E1:H1").PasteSpecial($xlPasteFormats)

$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# -- Index column (A2:A5) needs the same style as A2 already has -----------
$styleSource.Copy()
$q1.Range("A3:A5").PasteSpecial($xlPasteFormats)

$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# -- Data rows ---------------------------------------------------------------
# Fund codes & numeric-looking metrics are stored as TEXT in the source data
# (leading zeros in codes like "010709" must survive, and the percentages
# keep fixed decimal formatting) exactly like column B of the "2021-Q4"
# sheet. A leading apostrophe is the standard Excel way to force a
# numeric-looking value to stay text while keeping General formatting.
$q1.Range("B2").Value = "'010709"
$q1.Range("C2").Value = "安信医药健康主题股票A"
$q1.Range("D2").Value = "'20.94"
$q1.Range("E2").Value = "'89.03"
$q1.Range("F2").Value = "'3.14"
$q1.Range("G2").Value = "'0.6575"
$q1.Range("H2").Value = 8

$q1.Range("B3").Value = "'010710"
$q1.Range("C3").Value = "安信医药健康主题股票C"
$q1.Range("D3").Value = "'10.09"
$q1.Range("E3").Value = "'89.03"
$q1.Range("F3").Value = "'3.14"
$q1.Range("G3").Value = "'0.3168"
$q1.Range("H3").Value = 8

$q1.Range("B4").Value = "'009263"
$q1.Range("C4").Value = "华宝红利精选混合A"
$q1.Range("D4").Value = "'0.46"
$q1.Range("E4").Value = "'83.67"
$q1.Range("F4").Value = "'1.00"
$q1.Range("G4").Value = "'0.0046"
$q1.Range("H4").Value = 6

$q1.Range("B5").Value = "'010841"
$q1.Range("C5").Value = "华宝红利精选混合C"
$q1.Range("D5").Value = "'0.16"
$q1.Range("E5").Value = "'83.67"
$q1.Range("F5").Value = "'1.00"
$q1.Range("G5").Value = "'0.0016"
$q1.Range("H5").Value = 6

# ---------------------------------------------------------------------------
# 2) Add the new "总计" sheet right after "2022-Q1", rebuilding the totals
#    table with the new 2022-Q1 row on top of the old 2021-Q4 row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$styleSource.Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)
$styleSource.Copy()
$total.Range("A2:A3").PasteSpecial($xlPasteFormats)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.98

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 1.09
